$wb = $excel.ActiveWorkbook

# --- "2018 LEAVE CREDITS" sheet ---
$ws2018 = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws2018.Range("B2").Value = "DEMATERA, PEDRO JR. BAUSTISTA"
$ws2018.Range("C79").Value = 1.25

# --- "2017 LEAVE BALANCE" sheet ---
$ws2017 = $wb.Worksheets.Item("2017 LEAVE BALANCE")
$ws2017.Range("B20").Value = "SP(1-0-0)"
$ws2017.Range("K20").Value = "BDAY 3/27/2023"

$wb.Save()
